# Generate Report for Archive
# - Flip the "Ready for handoff" status label to "In Translation" everywhere it
#   is used (Overview summary columns + the per-locale Status column), and
#   shrink the now-narrower Status/summary columns to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Target rendered column width is 13.4101845877511 characters; this host's
# ColumnWidth setter snaps to the nearest 1/6-character increment, so feed it
# the character width (12.5) that lands on the closest achievable increment
# (13.333333333333334) to the real Excel-computed autofit width.
$newWidth  = 12.5

# --- Overview sheet: zh-cn / de-de summary columns (E & F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
